$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (pushes existing rows 4.. down to 5..),
# mirroring the new "Behav Proc in press" teaching entry that now sits
# between the "Assistant Professor" row and the old "Research Degree
# Project." row.
$ws.Rows(4).Insert()

# Row 2 ("Associate Professor" / Universidad El Bosque) now shows a
# concrete end year instead of the open-ended "2017 - Present" label.
$ws.Range("B2").Value = 2019

# Row 3 ("Assistant Professor") now spans 2017-2018 and references the
# "Quantitative Methods II" course (same text as row 2's E column).
$ws.Range("B3").Value = "2017-2018"
$ws.Range("E3").Value = "Quantitative Methods II (Psychology MSc)."

# New row 4: just the year and the course/description, no position /
# institution / location repeated (matches the blank A4,C4,D4 of the
# sibling sub-rows elsewhere in the sheet).
$ws.Range("A4").Clear()
$ws.Range("B4").Value = 2017
$ws.Range("C4").Clear()
$ws.Range("D4").Clear()
$ws.Range("E4").Value = "Quantitative Methods I (Psychology MSc)."

# Leave the selection where the author last left it.
$ws.Range("D20").Select() | Out-Null
